$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows for the candidates who withdrew from the race, from
# bottom to top so the earlier row numbers stay valid while deleting.
# Row 29: id 26 - Manuel Monsalve
# Row 20: id 17 - Francisco Huenchumilla
# Row 15: id 12 - David Fernandez
$ws.Rows.Item(29).EntireRow.Delete()
$ws.Rows.Item(20).EntireRow.Delete()
$ws.Rows.Item(15).EntireRow.Delete()

# Match the saved view state from the edit: scrolled down a bit with
# row 19 (now "Bernardo Javalquinto" -> ... ) selected as a whole row.
$ws.Application.Goto($ws.Range("A19"), $true)
$ws.Range("A19:XFD19").Select()
